$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: Status column (C), rows 2-4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column (C), rows 2-4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width changes (17.2159881591797 -> 13.4101848602295) ---
# Overview sheet: columns E and F
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
